$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
$d = $word.ActiveDocument

function Add-ParagraphXml($innerXml) {
    $lastPara = $d.Paragraphs.Last
    $r = $lastPara.Range
    [void]$r.InsertParagraphAfter()
    $newLast = $d.Paragraphs.Last
    $nr = $newLast.Range
    $nr.Collapse(0)
    [void]$nr.InsertXML("<w:p xmlns:w='$wNs'>$innerXml</w:p>")
}

# 1) paragraph with just a page break
Add-ParagraphXml "<w:r><w:br w:type=`"page`"/></w:r>"

# 2) "Total:" + 5 tabs, then "$136.5" in its own run
Add-ParagraphXml "<w:r><w:t>Total:</w:t><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/></w:r><w:r><w:t>`$136.5</w:t></w:r>"

# 3) "15% Discount:" + 4 tabs, then underlined "$(20.47)"
Add-ParagraphXml "<w:r><w:t>15% Discount:</w:t><w:tab/><w:tab/><w:tab/><w:tab/></w:r><w:r><w:rPr><w:u w:val=`"single`"/></w:rPr><w:t>`$(20.47)</w:t></w:r>"

# 4) bold "TOTAL:" + 4 tabs + "$116.03" all in one bold run
Add-ParagraphXml "<w:r><w:rPr><w:b/></w:rPr><w:t>TOTAL:</w:t><w:tab/><w:tab/><w:tab/><w:tab/><w:t>`$116.03</w:t></w:r>"
